$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the header cell in A1 (was the "fecha" column header) and clear
# all the date values that used to live under it in column A (rows 2-297).
# This also drops the now-unused "fecha" shared string, which shifts
# "Ins_Cupo" down one slot (from 701 to 700) - Excel keeps the H1 header
# cell pointing at the right text automatically because it is stored by
# value, not by shared-string index.
$ws.Range("A1:A297").ClearContents()

# Mirror the "select entire column A" selection state recorded in the
# sheet view after the edit.
$ws.Range("A1:A1048576").Select()
